# ---------------------------------------------------------------------------
# Applies the "Minutes23.03.2017.docx" edit described by the commit diff:
#   1. Merge the split "23/03" / "/2017" runs into a single "23/03/2017" run.
#   2. Split "Nazhoque" into "Naz" <space> <bookmark _GoBack> "hoque" and move
#      the _GoBack bookmark here (it used to sit at the end of the "Josh is
#      in the beginnings..." paragraph).
#   3. Append five new paragraphs after the "...teensy coding and
#      programming." paragraph with the meeting's outstanding action items.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "23/03" + "/2017" -> "23/03/2017"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("23/03/2017", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "23/03/2017", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Split "Nazhoque" and relocate the _GoBack bookmark.
# ---------------------------------------------------------------------------

# Remove the bookmark from its old position (end of the teensy paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate "Nazhoque" and figure out where to split it ("Naz" | "hoque").
$found = $d.Content
$found.Find.Execute("Nazhoque") | Out-Null
$nazStart = $found.Start
$splitPos = $nazStart + 3          # after "Naz"

# Insert the new _GoBack bookmark right between "Naz" and "hoque" - this
# naturally splits the single run into two runs ("Naz" / "hoque").
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Insert a space right after "Naz" (i.e. before the bookmark we just added).
$spaceRange = $d.Range($splitPos, $splitPos)
$spaceRange.InsertBefore(" ")

# Force the inserted space to live in its own run (otherwise it silently
# re-merges into the "Naz" run): toggle a character property on/off.
$spaceOnly = $d.Range($splitPos, $splitPos + 1)
$spaceOnly.Font.Bold = $true
$spaceOnly.Font.Bold = $false

# ---------------------------------------------------------------------------
# 3. Append the new paragraphs after "...teensy coding and programming. "
# ---------------------------------------------------------------------------
$found2 = $d.Content
$found2.Find.Execute("using a cartridge reader and the teensy coding and programming.") | Out-Null
$teensyPara = $found2.Paragraphs(1)

$teensyPara.Range.InsertParagraphAfter()

$newParaTexts = @(
    "The teensy code has been completed on time",
    "Task ",
    "The finalising of the piboy will be due the second week 13/04/2017",
    "And the finalisation of the walkthrough will be due on the 27/04/2017",
    "Naz is still to do the marketing documents."
)

$idx = $d.Paragraphs.Count
$d.Paragraphs($idx).Range.Text = $newParaTexts[0]

for ($i = 1; $i -lt $newParaTexts.Length; $i++) {
    $d.Paragraphs($idx).Range.InsertParagraphAfter()
    $idx = $d.Paragraphs.Count
    $d.Paragraphs($idx).Range.Text = $newParaTexts[$i]
}

Write-Output "edit complete"
